$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.701.40'
$ws.Range('E2').Value = '  -5.27%  '
$ws.Range('D3').Value = '''3.275.96'
$ws.Range('E3').Value = '  -6.51%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''176.89'
$ws.Range('E5').Value = '  -12.31%  '
$ws.Range('D6').Value = '''522.10'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('D7').Value = '''0.600'
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('D8').Value = '''3.266.48'
$ws.Range('E8').Value = '  -6.54%  '
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').Value = '''0.605'
$ws.Range('E10').Value = '  -7.62%  '
$ws.Range('D11').Value = '''57.43'
$ws.Range('E11').Value = '  -8.38%  '
$ws.Range('D12').Value = '''0.132'
$ws.Range('E12').Value = '  -7.96%  '
$ws.Range('D13').Value = '''0.0000256'
$ws.Range('E13').Value = '  -5.31%  '
$ws.Range('D14').Value = '''9.04'
$ws.Range('E14').Value = '  -7.89%  '
$ws.Range('D15').Value = '''3.801.73'
$ws.Range('E15').Value = '  -6.38%  '
$ws.Range('D16').Value = '''3.284.48'
$ws.Range('E16').Value = '  -6.14%  '
$ws.Range('E17').Value = '  -5.72%  '
$ws.Range('D18').Value = '''63.679.48'
$ws.Range('E18').Value = '  -4.89%  '
$ws.Range('D19').Value = '''17.33'
$ws.Range('E19').Value = '  -5.75%  '
$ws.Range('D20').Value = '''11.03'
$ws.Range('E20').Value = '  -6.66%  '
$ws.Range('D21').Value = '''0.950'
$ws.Range('E21').Value = '  -7.35%  '
$ws.Range('D22').Value = '''371.64'
$ws.Range('E22').Value = '  -5.10%  '
$ws.Range('D23').Value = '''3.75'
$ws.Range('E23').Value = '  -5.88%  '
$ws.Range('D24').Value = '''80.08'
$ws.Range('E24').Value = '  -3.28%  '
$ws.Range('D25').Value = '''11.00'
$ws.Range('E25').Value = '  -12.00%  '
$ws.Range('D26').Value = '''3.87'
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = '''2.66'
$ws.Range('E27').Value = '  -5.65%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''11.30'
$ws.Range('E28').Value = '  -8.01%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = '''8.29'
$ws.Range('E29').Value = '  -6.16%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''28.67'
$ws.Range('E30').Value = '  -7.52%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '''639.64'
$ws.Range('E31').Value = '  -5.79%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '''6.58'
$ws.Range('E32').Value = '  -5.79%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '''11.19'
$ws.Range('E33').Value = '  -4.59%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.105'
$ws.Range('E34').Value = '  -5.56%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''58.94'
$ws.Range('E35').Value = '  -7.48%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D37').Value = '''0.387'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '''36.35'
$ws.Range('E38').Value = '  -6.02%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '''0.0₃0693'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '''2.935.19'
$ws.Range('E41').Value = '  -4.46%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.124'
$ws.Range('E42').Value = '  -5.35%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '''2.44'
$ws.Range('E43').Value = '  -5.81%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = '''2.68'
$ws.Range('E44').Value = '  -10.32%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '''2.65'
$ws.Range('E45').Value = '  -4.57%  '
$ws.Range('D46').Value = '''0.0395'
$ws.Range('E46').Value = '  -1.18%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''3.03'
$ws.Range('E47').Value = '  +6.57%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '''2.78'
$ws.Range('E48').Value = '  +5.95%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '''0.125'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''135.29'
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '''2.44'
$ws.Range('E51').Value = '  -10.30%  '
